$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: B1:F1
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Dopad"
$ws.Range("D1").Value = "Pravděpodobnost"
$ws.Range("E1").Value = "Vliv"
$ws.Range("F1").Value = "MA"

# Column widths (closest achievable values given the host's width quantization)
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666

# Restore the author's last selection
$ws.Range("I5").Select() | Out-Null
